# Generate Report for Handback
#
# The underlying workbook stores repeated values (timestamps / priority codes)
# as shared strings. Several rows in the zh-cn / de-de / Overview sheets
# happened to contain identical text, so updating the shared value updates
# every cell that referenced it. Set all such cells explicitly so the
# resulting shared-string table matches.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# "Latest HO Xliff Generate Date" column G.
# G3 (8f50628f-752f-4517-bdd5-d2e09317b15b.md) previously shared the value
# "2016-09-01 18:17:25" with de-de!H3 and de-de!H5.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-01 18:18:35"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Priority column E: "ht" -> "mt" (shared by rows 3 and 5)
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
# Correspond Handoff Datetime column H: shared by rows 3 and 5
$wsZhCn.Range("H3").Value = "2016-09-01 18:18:29"
$wsZhCn.Range("H5").Value = "2016-09-01 18:18:29"
# Correspond Handback DateTime column K: shared by rows 3 and 5
$wsZhCn.Range("K3").Value = "2016-09-01 18:19:00"
$wsZhCn.Range("K5").Value = "2016-09-01 18:19:00"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# Priority column E: "ht" -> "mt" (shared by rows 3 and 5)
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
# Correspond Handoff Datetime column H: shared with Overview!G3 (rows 3 and 5)
$wsDeDe.Range("H3").Value = "2016-09-01 18:18:35"
$wsDeDe.Range("H5").Value = "2016-09-01 18:18:35"
# Correspond Handback DateTime column K: shared by rows 3 and 5
$wsDeDe.Range("K3").Value = "2016-09-01 18:19:19"
$wsDeDe.Range("K5").Value = "2016-09-01 18:19:19"
